# Daily "Updated symbol list" refresh for the cryptos sheet: prices /
# 1h-volume percentages are refreshed in place, and rows 41/42 swap their
# Coin/Link between BKEXToken and KickToken (with brand-new price/volume
# figures for each, not a simple value swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/1h) hold numeric- or percent-looking text
# that must remain plain text (matches original inlineStr cells with no
# explicit style). Force text entry via NumberFormat "@", then clear the
# format again so the cell keeps the default (unstyled) appearance.
function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue 2 4 "257.67"
Set-TextValue 2 5 "0.55%"

Set-TextValue 3 4 "27.10"
Set-TextValue 3 5 "-3.96%"

Set-TextValue 4 4 "4.898"
Set-TextValue 4 5 "-7.45%"

Set-TextValue 5 4 "0.05949"
Set-TextValue 5 5 "2.50%"

Set-TextValue 6 4 "6.688"
Set-TextValue 6 5 "-0.11%"

Set-TextValue 7 4 "0.8693"
Set-TextValue 7 5 "-0.28%"

Set-TextValue 8 4 "0.9602"
Set-TextValue 8 5 "5.94%"

Set-TextValue 9 4 "0.1413"
Set-TextValue 9 5 "0.15%"

Set-TextValue 10 4 "0.03552"
Set-TextValue 10 5 "3.74%"

Set-TextValue 11 4 "0.07178"
Set-TextValue 11 5 "-0.02%"

Set-TextValue 12 4 "0.03147"
Set-TextValue 12 5 "-0.36%"

Set-TextValue 13 4 "0.09252"
Set-TextValue 13 5 "0.15%"

Set-TextValue 14 4 "0.001542"
Set-TextValue 14 5 "0.26%"

Set-TextValue 15 4 "0.0006040"
Set-TextValue 15 5 "-94.26%"

Set-TextValue 16 4 "0.005983"
Set-TextValue 16 5 "1.59%"

Set-TextValue 17 5 "-0.61%"

Set-TextValue 18 4 "3.259"
Set-TextValue 18 5 "0.76%"

Set-TextValue 19 5 "-2.86%"

Set-TextValue 20 4 "0.3146"
Set-TextValue 20 5 "0.56%"

Set-TextValue 21 5 "-0.44%"

Set-TextValue 22 4 "3.538"
Set-TextValue 22 5 "0.44%"

Set-TextValue 23 4 "0.04273"
Set-TextValue 23 5 "2.61%"

Set-TextValue 24 5 "2.60%"

Set-TextValue 25 4 "0.001224"
Set-TextValue 25 5 "0.13%"

Set-TextValue 26 4 "0.004518"
Set-TextValue 26 5 "-9.16%"

Set-TextValue 27 5 "0.13%"

Set-TextValue 28 5 "-22.97%"

Set-TextValue 40 4 "0.03832"
Set-TextValue 40 5 "-1.01%"

$ws.Cells.Item(41, 2).Value = "KickToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue 41 4 "0.006580"
Set-TextValue 41 5 "14.27%"

$ws.Cells.Item(42, 2).Value = "BKEXToken"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue 42 4 "0.1104"
Set-TextValue 42 5 "0.66%"

Set-TextValue 43 4 "0.002200"
Set-TextValue 43 5 "-4.87%"

Set-TextValue 44 5 "0.00%"

Set-TextValue 45 4 "0.00005488"
Set-TextValue 45 5 "4.26%"

Set-TextValue 46 5 "0.13%"

Set-TextValue 47 5 "28.50%"

Set-TextValue 48 4 "0.002143"
Set-TextValue 48 5 "-1.74%"

Set-TextValue 49 5 "0.13%"

Set-TextValue 50 5 "0.13%"

